# CIDC-1278 first pass at redone docs
#
# Renames the ATACseq sheet/labels from "ATACSEQ" -> "ATACseq" casing,
# moves the active/selected tab from "Excluded Samples" to "Legend",
# updates each sheet's selection, and grows the two wrapped Legend rows
# that now need two lines of text.

$wb = $excel.ActiveWorkbook

$wsAnalysis = $wb.Worksheets.Item(1)   # "ATACSEQ Analysis" -> "ATACseq Analysis"
$wsExcluded = $wb.Worksheets.Item(2)   # "Excluded Samples"
$wsLegend   = $wb.Worksheets.Item(3)   # "Legend"
$wsDict     = $wb.Worksheets.Item(4)   # "Data Dictionary"

# --- Rename the main analysis tab (casing only) ------------------------
$wsAnalysis.Name = "ATACseq Analysis"

# --- Fix the "ATACSEQ" -> "ATACseq" casing in cell/legend text ---------
$wsAnalysis.Range("B7").Value = "ATACSeq Runs"

$wsLegend.Range("B2").Value = "Legend for tab 'ATACseq Analysis'"
$wsLegend.Range("B7").Value = "Section 'ATACSeq Runs' of tab 'ATACseq Analysis'"

# --- The two updated Legend rows now wrap onto two lines ----------------
$wsLegend.Rows.Item(2).RowHeight = 23.95
$wsLegend.Rows.Item(7).RowHeight = 23.95

# --- Update each sheet's selection (and, by selecting Legend last, make
#     it the active/visible tab instead of "Excluded Samples") ----------
[void]$wsAnalysis.Range("B7").Select()
[void]$wsExcluded.Range("B3").Select()
[void]$wsDict.Range("A1").Select()
[void]$wsLegend.Range("B10").Select()
